$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Quantity" column before the existing "Rank" column (old column E),
# shifting Rank/isActive/secType/expiration/strike/right/conID one column to the right.
$ws.Columns("E:E").Insert()

# New header + data for the Quantity column.
$ws.Range("E1").Value = "Quantity"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 10

# Updated trade data values.
$ws.Range("B3").Value = 111
$ws.Range("I2").Value = 20161021
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 215606586

# Give the trailing header cell (M1) the same bold header formatting as the rest of row 1.
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").Font.Size = 16

# Widen the "expiration" column (now column I) to fit its contents.
$ws.Columns("I:I").ColumnWidth = 11.71

# Move the active selection.
$ws.Range("F4").Select()
